# Populate new "phenotype" (K) and "genotype" (L) status columns on the
# "clinical values" sheet (3rd sheet of the workbook).
#
# NOTE: the order in which cells are written below is deliberately chosen
# (matches how the values were originally entered) so that new shared
# strings are appended to xl/sharedStrings.xml in the same order as the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Header row
$ws.Range("K1").Value = "phenotype"
$ws.Range("L1").Value = "genotype"

# First data row
$ws.Range("K2").Value = "phenotypic-manifestation"
$ws.Range("L2").Value = "mutation-carrier"

# Genotype (L) column for the remaining rows
$ws.Range("L3").Value  = "no-mutation-carrier"
$ws.Range("L4").Value  = "not-genotyped"
$ws.Range("L5").Value  = "no-mutation-carrier"
$ws.Range("L6").Value  = "not-genotyped"
$ws.Range("L7").Value  = "mutation-carrier"
$ws.Range("L8").Value  = "no-mutation-carrier"
$ws.Range("L9").Value  = "not-genotyped"
$ws.Range("L10").Value = "mutation-carrier"
$ws.Range("L11").Value = "not-genotyped"
$ws.Range("L12").Value = "mutation-carrier"

# Phenotype (K) column for the remaining rows
$ws.Range("K3").Value  = "no-phenotypic-manifestation"
$ws.Range("K4").Value  = "sudden-cardiac-death"
$ws.Range("K8").Value  = "death-from-other-causes"
$ws.Range("K5").Value  = "ressucitated-scd-or-vf"
$ws.Range("K6").Value  = "phenotypic-manifestation"
$ws.Range("K7").Value  = "phenotypic-manifestation"
$ws.Range("K9").Value  = "phenotypic-manifestation"
$ws.Range("K10").Value = "ressucitated-scd-or-vf"
$ws.Range("K11").Value = "phenotypic-manifestation"
$ws.Range("K12").Value = "sudden-cardiac-death"

# Column widths for the two new columns (closest achievable values to the
# target OOXML widths of 23.6640625 and 18, given the runtime's internal
# pixel-based quantization of ColumnWidth).
$ws.Columns.Item(11).ColumnWidth = 22.83
$ws.Columns.Item(12).ColumnWidth = 17.166666666666668

# Move / update the active selection to K10, matching the saved view state.
$ws.Range("K10").Select()
